$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, pushing existing rows 16-27 down to 17-28.
$ws.Rows.Item(16).Insert()

# Populate the new row 16 with the latest week's data.
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 44413
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112035
$ws.Range("G16").Value = "Bruselas (repollito)"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 40
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 25000
$ws.Range("N16").Value = "$/malla 10 kilos"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 2500
$ws.Range("Q16").Value = 10
$ws.Range("R16").Value = "Hortaliza"
